# Diary workbook update: add a new diary entry (27 loka) as row 16,
# plus a trailing formatted-but-empty row 17, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New entry: row 16 ---------------------------------------------------
# Pull cell formatting (number formats / wrap text) from row 6, which is an
# existing entry with the same shape (A..F populated, same row height).
$ws.Range("B6:F6").Copy() | Out-Null
$ws.Range("B16:F16").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Values are assigned in the same order the shared-string table picks them
# up (date, content, time, notes, meta, quality) so the new <si> entries
# land at the expected indices.
$ws.Range("A16").Value = "27 loka"
$ws.Range("C16").Value = "Fireworkit näkymään modernimmalla c++, pikakurssi C++ iteraattoreihin tekoälyn johdolla, pikakertaus OpenGL primitiiveihin"
$ws.Range("B16").Value = "12.45-15.15, 16.15-17.15, 18.15-20.15, 20.30-21.00"
$ws.Range("E16").Value = "Koodistakin alkaa löytyä vähän järkeä, kun sai tuon uniformin käytettyä uudestaan firework scenessä. Nyt on oikeasti hyvä fiilis jatkaa eteenpäin, kun sai kaiken toimimaan ja hieman modernimmalla c++:lla. Aika hakusessahan tuo vielä on, ja virheen löytäminen on työn ja tuskan takana."
$ws.Range("F16").Value = "Tästä jatketaan kohti kytkettyjä kappaleita!"
$ws.Range("D16").Value = "No, siellä 80-90% mentiin, ei voi ymmärtää mitä ei ymmärrä ja tässä haasteena on modernisoida koodikantaa mitä ei kunnolla ymmärrä, eikä sitä modernisointiakaan vielä hanskaa. Taistelu kerrallaan toivon mukaan kehittyy myös tässä"
$ws.Range("G16").Value = 6

$ws.Rows.Item(16).RowHeight = 116

# --- Trailing formatted row 17 -------------------------------------------
# B17 carries the same style as B11 (time-format cell), left empty.
$ws.Range("B11").Copy() | Out-Null
$ws.Range("B17").PasteSpecial(-4122) | Out-Null       # xlPasteFormats

# --- View state ------------------------------------------------------------
# (topLeftCell scrolling isn't exposed by this host's Window object, so only
# the active selection is reproducible here.)
$ws.Range("A9").Select() | Out-Null
$ws.Range("D17").Select() | Out-Null
